# "Added the missing runs"
#
# The original sheet listed the "Rise" runs for h1/f1200 in pressure order,
# but the p14 run (row 26: R_h1_f1200_1_p14) was sitting where a "p10" run
# should have been. This edit:
#   1. Removes row 26 (R_h1_f1200_1_p14), which shifts the rows below it
#      (27-39) up by one (26-38) - this also folds in a data correction on
#      what is now row 31 (R_h2_f1200_1_p13 Quality: Medium -> Good).
#   2. Appends three rows at the bottom of the table:
#        - row 39: a new "Later" section marker row for F_h2_f1000_1_s
#        - row 40: the previously-missing R_h1_f750_1_p10 run
#        - row 41: the R_h1_f1200_1_p14 run that used to be row 26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove old row 26 (R_h1_f1200_1_p14), shifting 27-39 up to 26-38 ---
$ws.Rows("26:26").Delete()

# --- Step 2: append the new rows at the bottom (39, 40, 41) ---
# (B40's string is interned into the shared-string table first so it lands
# at the same index the source workbook uses: R_h1_f750_1_p10=49, Later=50.)

# Row 40: previously-missing R_h1_f750_1_p10 run
$ws.Range("B40").Value = "R_h1_f750_1_p10"
$ws.Range("C40").HorizontalAlignment = -4108
$ws.Range("C40").VerticalAlignment = -4108
$ws.Range("D40").HorizontalAlignment = -4108
$ws.Range("D40").VerticalAlignment = -4108
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = "Good"

# Row 39: new "Later" marker row
$ws.Range("A39").Value = "Later"
$ws.Range("B39").Value = "F_h2_f1000_1_s"
$ws.Range("C39").HorizontalAlignment = -4108
$ws.Range("C39").VerticalAlignment = -4108
$ws.Range("D39").HorizontalAlignment = -4108
$ws.Range("D39").VerticalAlignment = -4108
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = "Good"

# Row 41: R_h1_f1200_1_p14 run moved here from the old row 26
$ws.Range("B41").Value = "R_h1_f1200_1_p14"
$ws.Range("C41").HorizontalAlignment = -4108
$ws.Range("C41").VerticalAlignment = -4108
$ws.Range("D41").HorizontalAlignment = -4108
$ws.Range("D41").VerticalAlignment = -4108
$ws.Range("C41").Value = 4
$ws.Range("D41").Value = 4
$ws.Range("E41").Value = "Good"

# --- View state: match the author's final selection/scroll ---
[void]$ws.Range("G28").Select()
